$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the width of the column immediately to the left (J), so the
# newly inserted column inherits the same width (matches Excel's default
# insert-column behaviour of carrying over the left neighbour's formatting).
$existingWidth = $ws.Columns.Item(10).ColumnWidth

# Insert a new column before column K (11th column) to hold the new
# "VAT Partner" field, shifting the existing columns (K..R) right to (L..S)
$ws.Columns.Item(11).Insert()

$ws.Columns.Item(11).ColumnWidth = $existingWidth

$ws.Cells.Item(1, 11).Value = "VAT Partner"

$ws.Range("K2").Select()

# Update the filter database defined name to extend over the new column
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$S`$1"

